# Added feature to add a product and write entry in excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProducts = @(
    @{ Id = 6; Name = "handwash"; Brand = "kenvue";   Price = 1000; Qty = 5  },
    @{ Id = 7; Name = "lotion";   Brand = "jnj";      Price = 1000; Qty = 10 },
    @{ Id = 8; Name = "skincare"; Brand = "kolly";    Price = 100;  Qty = 12 },
    @{ Id = 9; Name = "wizard";   Brand = "harry";    Price = 100;  Qty = 12 }
)

foreach ($product in $newProducts) {
    $row = $product.Id
    $ws.Cells.Item($row, 1).Value = $product.Id
    $ws.Cells.Item($row, 2).Value = $product.Name
    $ws.Cells.Item($row, 3).Value = $product.Brand
    $ws.Cells.Item($row, 4).Value = $product.Price
    $ws.Cells.Item($row, 5).Value = $product.Qty
}
